# Hybrid bold + color highlighting for quantitative impact metrics
# (percentages, dollar amounts, large numbers) across achievements and
# work-experience bullet points, mirroring the DOCX bold+color treatment.

$d = $word.ActiveDocument
$highlightColor = 5258796   # BGR packed value for RGB(0x2C,0x3E,0x50) -> w:color "2C3E50"

function Find-ParagraphRange($doc, $includeText, $excludeText) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $t = $p.Range.Text
        if ($t -notlike ("*" + $includeText + "*")) {
            continue
        }
        if ($excludeText -and ($t -like ("*" + $excludeText + "*"))) {
            continue
        }
        return $p.Range
    }
    return $null
}

function Apply-MetricHighlights($doc, $paraRange, $segments) {
    if ($paraRange -eq $null) {
        Write-Output "PARAGRAPH NOT FOUND"
        return
    }
    $paraEnd = $paraRange.End
    $pos = $paraRange.Start
    foreach ($seg in $segments) {
        $r = $doc.Range($pos, $paraEnd)
        $found = $r.Find.Execute($seg, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) {
            Write-Output ("NOT FOUND: " + $seg)
            continue
        }
        $r.Font.Bold = $true
        $r.Font.Color = $highlightColor
        $pos = $r.End
    }
}

# 1) "Discovered systematic race coding errors ... from 23% to 64%"
$p1 = Find-ParagraphRange $d 'Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning' $null
Apply-MetricHighlights $d $p1 @('23%', '64%')

# 2) "Achieved 87% prediction accuracy ... margins from ±4.2% to ±2.1%"
$p2 = Find-ParagraphRange $d 'Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins' $null
Apply-MetricHighlights $d $p2 @('87%', '71%', '±4.2%', '±2.1%')

# 3) "Wrote RFP and analyzed bids from 1,200 vendors ..."
$p3 = Find-ParagraphRange $d 'Wrote RFP and analyzed bids from 1,200 vendors' $null
Apply-MetricHighlights $d $p3 @('1,200')

# 4) "... became the $400M Polling Consortium Database ... valued at $1B+"
$p4 = Find-ParagraphRange $d 'Created comprehensive meta-analysis framework handling millions of survey responses' $null
Apply-MetricHighlights $d $p4 @('$400M', '$1B')

# 5) "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
$p5 = Find-ParagraphRange $d 'Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations' $null
Apply-MetricHighlights $d $p5 @('73.5%', '$4.7M')

# 6) "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" (short form, no trailing clause)
$p6 = Find-ParagraphRange $d 'Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%' 'reducing polling'
Apply-MetricHighlights $d $p6 @('87%', '71%')

Write-Output "Quantitative metrics highlighting applied."
